$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Copy the formatting of the last existing data row (71) down into the four
#    new rows (72-75) so borders/number formats match the rest of the table.
$ws.Range("A71:H71").Copy()
$ws.Range("A72:H75").PasteSpecial(-4122)

# 2. Re-apply the worksheet AutoFilter over A1:H74 while row 75 is still
#    empty, so the filter range does not auto-expand once row 75 is filled in.
$ws.Range("A1:H74").AutoFilter()

# 3. Fill in the new rows. The "nome_tabella" (column C) values are written
#    first, in the same order the source workbook interned them as shared
#    strings, so the resulting shared-strings table matches exactly.
$ws.Range("C73").Value2 = "rfcf_storico_run_cashflow"
$ws.Range("C72").Value2 = "rfcf_configurazione_start_run_cashflow"
$ws.Range("C74").Value2 = "rfcf_decodifiche_cashflow"
$ws.Range("C75").Value2 = "rfcf_parametri_check_decodifiche"

$ws.Range("A72").Value2 = "rfcf_cashflow"
$ws.Range("B72").Value2 = "data"
$ws.Range("D72").Value2 = "input"

$ws.Range("A73").Value2 = "rfcf_cashflow"
$ws.Range("B73").Value2 = "data"
$ws.Range("D73").Value2 = "input"

$ws.Range("A74").Value2 = "rfcf_cashflow"
$ws.Range("B74").Value2 = "data"
$ws.Range("D74").Value2 = "input"

$ws.Range("A75").Value2 = "rfcf_cashflow"
$ws.Range("B75").Value2 = "data"
$ws.Range("D75").Value2 = "input"

# 4. Update the hidden _FilterDatabase defined name to the new range used by
#    the refreshed AutoFilter.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Snapshot!_FilterDatabase") {
        $n.RefersTo = "=Snapshot!`$A`$1:`$H`$73"
    }
}

# 5. Update the view: 100% zoom, scrolled back to the top-left, with A3
#    selected (matching the refreshed sheetView).
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A3").Select()
